$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) column C holds a date serial number for every
# data row (rows 2-321). The workbook was refreshed a day later, so every
# value moves from 45179 (2023-09-10) to 45180 (2023-09-11).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 321 }

$rng = $ws.Range("C2:C$lastRow")
$rng.Value = 45180
